# This script shifts the existing sensor readings (columns C:H, rows 2-21)
# down by 4 rows, then populates the newly freed rows 2-5 with 4 new
# "falling" sample readings, matching the target dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current values of C2:H17 (16 rows) before overwriting anything,
# since they will be moved down to C6:H21.
$srcRange = $ws.Range("C2:H17")
$srcValues = $srcRange.Value()

# Write those captured values into C6:H21 (shift down by 4 rows).
$dstRange = $ws.Range("C6:H21")
$dstRange.Value = $srcValues

# New readings for the 4 freshly inserted rows (2-5).
$newRows = @(
    @(-3.012916564941406, 8.089370727539062, -0.1633265316486358, 0.03629761248826986,  0.01907121278345579,  0.05546045627444995),
    @(-3.395848751068115, 8.023316383361816,  0.0382503271102905, 0.05165476366877556, -0.0003787364251911958, 0.03377473920583711),
    @(-3.384797096252441, 7.934267520904541,  0.07479587197303771, 0.02702467799186697, -0.02729956846684218, -0.006963863894343374),
    @(-3.632324695587158, 7.965863227844238,  0.0220168232917785, 0.02345722466707222, -0.01078177168965329, -0.01979203335940831)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = 2 + $i
    $values = $newRows[$i]
    for ($col = 0; $col -lt $values.Count; $col++) {
        # Column C is index 3 (A=1, B=2, C=3, ...)
        $ws.Cells.Item($row, 3 + $col).Value = $values[$col]
    }
}
